$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two dollar totals affected by the new transaction
$ws.Range("B4").Value = 401
$ws.Range("B6").Value = 399

# Add the new transaction log entry in row 14, matching the formatting
# of the existing log rows (copy format from A13, the row above it)
$ws.Range("A14").Value = "22.02.2025 - Out of PO Vassell karşılığında MaltaSpor'a 1 Dolar vermiştir. (401-399)"
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The saved file no longer flags column A's width as auto ("best fit") -
# nudge the column width so it is stored as an explicit custom width.
$ws.Columns.Item(1).ColumnWidth = 87.25

# Update the saved cell selection
[void]$ws.Range("A18").Select()
